$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update test data values ---
# Order matters: new shared-string entries are appended to the shared
# strings table in the order they are first written, so we write them in
# the same order they appear in the target workbook (Company, Address,
# AddrAlias for row2, AddrAlias for row3, then the two DateOfBirth values)
# so the shared string table layout matches the target file.
$ws.Range("K2").Value = "MNC company"
$ws.Range("L2").Value = "Outer Ring Road"
$ws.Range("U2").Value = "Indian,India"
$ws.Range("U3").Value = "Indian ,India,Bangalore"
$ws.Range("F2").Value = "12-03-2o11"
# F3 already carries a text/quote-prefix format (it shows a value that looks
# like a date). Prefix with an apostrophe so Excel keeps storing it as plain
# text using the existing style instead of re-parsing it as a date serial.
$ws.Range("F3").Value = "'01-05-2014"

# --- Update view / selection ---
$ws.Range("D3").Select()

# --- Update column widths ---
$ws.Columns.Item(11).ColumnWidth = 19.15
$ws.Columns.Item(13).ColumnWidth = 10.855
